# Update the "Förändrad" (Changed) date column (C) from 45563 to 45564
# for all data rows (2 through 29) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45563) {
        $cell.Value = 45564
    }
}
